$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 — reuse G1's formatting (bold, centered,
# bordered header style) by copy/paste-special of formats only, so the
# existing style index is reused rather than a near-duplicate being minted.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column values for each data row (era data update)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
